# Rename a handful of country labels in column A to match the naming
# convention used by the new case-map dataset (column B country codes are
# unchanged; only the display names in column A are updated, in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renames = @{
    16  = "Bahamas"
    47  = "Congo (Kinshasa)"
    48  = "Congo (Brazzaville)"
    56  = "Czechia"
    75  = "Gambia"
    193 = "Eswatini"
    197 = "Taiwan*"
    213 = "US"
    220 = "West Bank and Gaza"
}

foreach ($row in $renames.Keys) {
    $ws.Cells.Item($row, 1).Value = $renames[$row]
}
